# Updates the cryptocurrency price/volume table (and one rebranded coin row)
# on Sheet1, mirroring the scraped-data refresh from the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price"/"Volume(1h)" columns (and the B51:C51 coin name/link pair we
# are about to overwrite) are stored as plain text in the workbook, even
# though many of the values look numeric (e.g. "0.999", "1.00"). Excel's
# COM Value setter auto-converts number-looking strings to real numbers,
# which would corrupt values like trailing zeros ("1.00" -> 1) and the
# dotted-thousands prices ("60.692.94"). Temporarily force the range to
# Text format so every assignment below is kept as a literal string, then
# clear the formatting back to the workbook default afterwards so no stray
# number-format style is left behind.
$textRange = $ws.Range("B2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "60.692.94"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "2.622.91"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "569.85"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").Value = "142.24"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.599"
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("D9").Value = "2.622.41"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").Value = "6.54"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").Value = "0.149"
$ws.Range("E13").Value = "  -6.80%  "
$ws.Range("D14").Value = "3.084.28"
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "60.654.03"
$ws.Range("E15").Value = "  +2.42%  "
$ws.Range("D16").Value = "23.33"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("D18").Value = "2.613.03"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "11.36"
$ws.Range("E19").Value = "  +9.60%  "
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").Value = "346.56"
$ws.Range("E21").Value = "  +2.84%  "
$ws.Range("D22").Value = "6.97"
$ws.Range("E22").Value = "  +8.33%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "0.529"
$ws.Range("E24").Value = "  +13.15%  "
$ws.Range("D25").Value = "63.39"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "0.987"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  +4.53%  "
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("E30").Value = "  +9.12%  "
$ws.Range("D31").Value = "6.39"
$ws.Range("E31").Value = "  +4.14%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("E35").Value = "  +4.83%  "
$ws.Range("D36").Value = "0.971"
$ws.Range("E36").Value = "  +10.55%  "
$ws.Range("E37").Value = "  +4.21%  "
$ws.Range("E38").Value = "  +8.81%  "
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("E40").Value = "  +3.67%  "
$ws.Range("D41").Value = "0.852"
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("D42").Value = "295.16"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("D43").Value = "138.82"
$ws.Range("E43").Value = "  +4.77%  "
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D46").Value = "0.606"
$ws.Range("D47").Value = "19.74"
$ws.Range("E47").Value = "  +3.82%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  +2.94%  "
$ws.Range("D50").Value = "19.89"
$ws.Range("E50").Value = "  +6.73%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "10.72"
$ws.Range("E51").Value = "  +0.62%  "

$textRange.ClearFormats()
